$wb = $excel.ActiveWorkbook

# Existing sheets: "总计" (Total), "2022-Q3" (fund detail data), "2021-Q3" (fund detail data)
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# --- Step 1: split off the current "2022-Q3" fund-detail sheet ---
# Duplicate it first so the existing (now historical) values survive unchanged
# as the new "2022-Q3" sheet, placed right after the sheet that will become "2022-Q4".
$wsQ3.Copy([System.Reflection.Missing]::Value, $wsQ3)
$wsNewQ3 = $wb.ActiveSheet
$wsNewQ3.Name = "2022-Q3 temp"

# Turn the original sheet into "2022-Q4" and refresh it with the latest figures.
# (Re-fetch the sheet reference fresh by name since the structural Copy() above can
# invalidate previously captured worksheet references.)
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Name = "2022-Q4"

$wsQ3.Range("D2").NumberFormat = "@"
$wsQ3.Range("D2").Value = "50.81"
$wsQ3.Range("D2").Style = "Normal"

$wsQ3.Range("E2").NumberFormat = "@"
$wsQ3.Range("E2").Value = "94.73"
$wsQ3.Range("E2").Style = "Normal"

$wsQ3.Range("F2").NumberFormat = "@"
$wsQ3.Range("F2").Value = "5.26"
$wsQ3.Range("F2").Style = "Normal"

$wsQ3.Range("G2").NumberFormat = "@"
$wsQ3.Range("G2").Value = "2.6726"
$wsQ3.Range("G2").Style = "Normal"

$wsQ3.Range("H2").Value = 9

# Give the duplicate (holding the old values) its final name.
$wsNewQ3 = $wb.Worksheets.Item("2022-Q3 temp")
$wsNewQ3.Name = "2022-Q3"

# --- Step 2: update the "总计" (Total) summary sheet ---
$wsTotal = $wb.Worksheets.Item("总计")

# Insert a new row above the existing 2022-Q3 summary row, shifting the rest down.
$wsTotal.Range("A3:D3").Insert()

# Copy the old row2 (2022-Q3 summary figures) down into the newly inserted row3.
$wsTotal.Range("B2").Copy()
$wsTotal.Range("B3").PasteSpecial(-4163)
$wsTotal.Range("C2").Copy()
$wsTotal.Range("C3").PasteSpecial(-4163)
$wsTotal.Range("D2").Copy()
$wsTotal.Range("D3").PasteSpecial(-4163)
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Overwrite row2 in place with the new 2022-Q4 figures.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 2.67

# Fix up the running index column (0, 1, 2).
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2

# Restore the originally active/selected sheet.
$wb.Worksheets.Item("2021-Q3").Activate()
